$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 871 (2026/12/29 火 block),
# shifting the existing rows 871-912 down to 873-914.
$ws.Range("871:872").EntireRow.Insert()

# Populate the two newly inserted rows with the new data points for
# 2026/02/27 (金). Force column A to stay plain text (it otherwise gets
# auto-converted to a date serial number), then reset the style so the
# cell matches the plain (unstyled) look of its sibling date cells.
$ws.Range("A871").NumberFormat = "@"
$ws.Range("A871").Value = "2026/02/27"
$ws.Range("A871").Style = "Normal"
$ws.Range("B871").Value = "金"
$ws.Range("C871").Value = 19
$ws.Range("D871").Value = 24

$ws.Range("A872").NumberFormat = "@"
$ws.Range("A872").Value = "2026/02/27"
$ws.Range("A872").Style = "Normal"
$ws.Range("B872").Value = "金"
$ws.Range("C872").Value = 22
$ws.Range("D872").Value = 24
